# Updates cryptos list prices/volumes per the GitHub Actions scraping run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value would otherwise be auto-parsed by Excel
# as a plain number (losing the original decimal-string formatting / introducing
# floating point rounding). Force these to Text before assigning so the stored
# value matches the scraped display string exactly.
$textRows = @(4,5,6,7,8,9,10,11,12,13,15,16,17,18,19,20,21,22,24,25,27,28,29,30,31,32,33,34,35,36,37,38,39,41,42,44,45,46,47,48,49,50,51)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.265.66"
$ws.Range("E2").Value = "  +5.52%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.911.07"
$ws.Range("E3").Value = "  +1.99%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.58%  "

# Row 5 - BNB
$ws.Range("D5").Value = "329.03"
$ws.Range("E5").Value = "  +4.60%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.56%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.5186"
$ws.Range("E7").Value = "  +2.12%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.4040"
$ws.Range("E8").Value = "  +3.46%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.08497"
$ws.Range("E9").Value = "  +1.64%  "

# Row 10 - was OKB, now Polygon (rows 10/11 swapped order)
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.125"
$ws.Range("E10").Value = "  +1.78%  "

# Row 11 - was Polygon, now OKB
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "42.79"
$ws.Range("E11").Value = "  +1.24%  "

# Row 12 - Solana
$ws.Range("D12").Value = "23.67"
$ws.Range("E12").Value = "  +16.39%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "6.434"
$ws.Range("E13").Value = "  +4.20%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.914.88"
$ws.Range("E14").Value = "  +2.29%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "7.386"
$ws.Range("E15").Value = "  +1.98%  "

# Row 16 - BinanceUSD
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.71%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "95.15"
$ws.Range("E17").Value = "  +2.19%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.00001115"
$ws.Range("E18").Value = "  +1.54%  "

# Row 19 - TRON
$ws.Range("D19").Value = "0.06717"
$ws.Range("E19").Value = "  +0.16%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "18.47"
$ws.Range("E20").Value = "  +4.91%  "

# Row 21 - Dai
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.50%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.016"
$ws.Range("E22").Value = "  +1.65%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "30.274.41"

# Row 24 - Cosmos
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  +2.45%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.221"
$ws.Range("E25").Value = "  -0.33%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").Value = "2.140.32"
$ws.Range("E26").Value = "  +2.59%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "21.95"
$ws.Range("E27").Value = "  +6.67%  "

# Row 28 - Monero
$ws.Range("D28").Value = "162.50"
$ws.Range("E28").Value = "  +3.26%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "2.403"
$ws.Range("E29").Value = "  -0.20%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "129.67"
$ws.Range("E30").Value = "  +2.83%  "

# Row 31 - ImmutableX
$ws.Range("D31").Value = "1.107"
$ws.Range("E31").Value = "  +5.99%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "0.1064"

# Row 33 - Filecoin
$ws.Range("D33").Value = "6.015"
$ws.Range("E33").Value = "  +4.04%  "

# Row 34 - HuobiToken
$ws.Range("D34").Value = "3.645"
$ws.Range("E34").Value = "  +0.24%  "

# Row 35 - VeChain
$ws.Range("D35").Value = "0.02494"
$ws.Range("E35").Value = "  +1.92%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "0.06587"
$ws.Range("E36").Value = "  +0.88%  "

# Row 37 - Algorand
$ws.Range("D37").Value = "0.2214"
$ws.Range("E37").Value = "  +2.42%  "

# Row 38 - was InternetComputer(DFINITY), now ARBITRUM (rows 38/39 swapped order)
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.234"
$ws.Range("E38").Value = "  +3.28%  "

# Row 39 - was ARBITRUM, now InternetComputer(DFINITY)
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "1.234"
$ws.Range("E39").Value = "  +2.87%  "

# Row 40 - Aptos
$ws.Range("E40").Value = "  +7.52%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "8.807"
$ws.Range("E41").Value = "  -2.15%  "

# Row 42 - TheSandbox
$ws.Range("D42").Value = "0.6533"
$ws.Range("E42").Value = "  +2.47%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  -0.01%  "

# Row 44 - Decentraland
$ws.Range("D44").Value = "0.6145"
$ws.Range("E44").Value = "  +2.70%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "13.31"
$ws.Range("E45").Value = "  +2.44%  "

# Row 46 - PancakeSwap
$ws.Range("D46").Value = "3.737"
$ws.Range("E46").Value = "  +1.78%  "

# Row 47 - NEARProtocol
$ws.Range("D47").Value = "2.072"
$ws.Range("E47").Value = "  +3.48%  "

# Row 48 - EOS
$ws.Range("D48").Value = "1.246"
$ws.Range("E48").Value = "  +1.92%  "

# Row 49 - Quant
$ws.Range("D49").Value = "125.42"
$ws.Range("E49").Value = "  +2.82%  "

# Row 50 - WEMIXTOKEN
$ws.Range("D50").Value = "1.160"
$ws.Range("E50").Value = "  -1.65%  "

# Row 51 - Aave
$ws.Range("D51").Value = "79.53"
$ws.Range("E51").Value = "  +4.37%  "
